$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.090.12"
$ws.Range("E2").Value = "  -2.48%  "
$ws.Range("D3").Value = "3.600.63"
$ws.Range("E3").Value = "  -2.58%  "
$ws.Range("E4").Value = "  +0.27%  "
$ws.Range("D5").Value = "'582.36"
$ws.Range("E5").Value = "  -2.58%  "
$ws.Range("D6").Value = "'179.26"
$ws.Range("E6").Value = "  -2.70%  "
$ws.Range("D7").Value = "'0.608"
$ws.Range("E7").Value = "  -3.50%  "
$ws.Range("E8").Value = "  -0.14%  "
$ws.Range("D9").Value = "'0.669"
$ws.Range("E9").Value = "  -7.03%  "
$ws.Range("D10").Value = "'0.144"
$ws.Range("E10").Value = "  -11.22%  "
$ws.Range("D11").Value = "'52.99"
$ws.Range("E11").Value = "  -6.13%  "
$ws.Range("D12").Value = "'0.0000251"
$ws.Range("E12").Value = "  -13.88%  "
$ws.Range("D13").Value = "'9.87"
$ws.Range("E13").Value = "  -7.63%  "
$ws.Range("D14").Value = "4.198.42"
$ws.Range("E14").Value = "  -2.65%  "
$ws.Range("D15").Value = "3.619.95"
$ws.Range("E15").Value = "  -3.11%  "
$ws.Range("D16").Value = "'0.125"
$ws.Range("E16").Value = "  -0.40%  "
$ws.Range("D17").Value = "67.017.78"
$ws.Range("E17").Value = "  -2.50%  "
$ws.Range("D18").Value = "'18.23"
$ws.Range("E18").Value = "  -5.82%  "
$ws.Range("D19").Value = "'12.14"
$ws.Range("E19").Value = "  -5.53%  "
$ws.Range("E20").Value = "  -5.52%  "
$ws.Range("D21").Value = "'389.79"
$ws.Range("E21").Value = "  -4.96%  "
$ws.Range("D22").Value = "'4.27"
$ws.Range("D23").Value = "'84.94"
$ws.Range("E23").Value = "  -4.39%  "
$ws.Range("D24").Value = "'2.84"
$ws.Range("E24").Value = "  -6.55%  "
$ws.Range("B25").Value = "LEO"
$ws.Range("C25").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D25").Value = "'6.06"
$ws.Range("E25").Value = "  -0.36%  "
$ws.Range("B26").Value = "InternetComputer(DFINITY)"
$ws.Range("C26").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D26").Value = "'12.17"
$ws.Range("E26").Value = "  -5.01%  "
$ws.Range("D27").Value = "'10.22"
$ws.Range("E27").Value = "  -6.26%  "
$ws.Range("D28").Value = "'3.58"
$ws.Range("E28").Value = "  -13.06%  "
$ws.Range("D29").Value = "'8.92"
$ws.Range("E29").Value = "  -5.83%  "
$ws.Range("D30").Value = "'30.94"
$ws.Range("E30").Value = "  -5.62%  "
$ws.Range("D31").Value = "'6.67"
$ws.Range("E31").Value = "  -8.46%  "
$ws.Range("D32").Value = "'66.50"
$ws.Range("E32").Value = "  +2.10%  "
$ws.Range("D33").Value = "'11.80"
$ws.Range("E33").Value = "  -5.13%  "
$ws.Range("D34").Value = "'0.111"
$ws.Range("E34").Value = "  -5.82%  "
$ws.Range("D35").Value = "'581.58"
$ws.Range("E35").Value = "  -4.41%  "
$ws.Range("D36").Value = "'41.07"
$ws.Range("E36").Value = "  -5.88%  "
$ws.Range("D37").Value = "'1.00"
$ws.Range("E37").Value = "  +0.13%  "
$ws.Range("E38").Value = "  -0.20%  "
$ws.Range("D39").Value = "'0.372"
$ws.Range("E39").Value = "  -7.37%  "
$ws.Range("D40").Value = "0.0₃0729"
$ws.Range("E40").Value = "  -18.92%  "
$ws.Range("D41").Value = "'0.131"
$ws.Range("E41").Value = "  -4.51%  "
$ws.Range("D42").Value = "'2.75"
$ws.Range("E42").Value = "  -9.59%  "
$ws.Range("D43").Value = "'0.0408"
$ws.Range("E43").Value = "  -7.20%  "
$ws.Range("D44").Value = "2.682.00"
$ws.Range("E44").Value = "  -3.34%  "
$ws.Range("B45").Value = "ApeXProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D45").Value = "'3.07"
$ws.Range("E45").Value = "  -3.95%  "
$ws.Range("B46").Value = "Fetch.AI"
$ws.Range("C46").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D46").Value = "'2.37"
$ws.Range("E46").Value = "  -13.05%  "
$ws.Range("B47").Value = "Stellar"
$ws.Range("C47").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D47").Value = "'0.128"
$ws.Range("E47").Value = "  -4.68%  "
$ws.Range("D48").Value = "'2.54"
$ws.Range("E48").Value = "  -6.81%  "
$ws.Range("D49").Value = "'135.68"
$ws.Range("E49").Value = "  -4.30%  "
$ws.Range("D50").Value = "'8.13"
$ws.Range("E50").Value = "  -11.93%  "
$ws.Range("D51").Value = "'2.56"
$ws.Range("E51").Value = "  -7.17%  "
